# Fix: Elective lecture and tutorial scheduling
# Applies to the six timetable sheets (Regular/PreMid/PostMid x Section A/B):
#   - CS303 (Lab) room change: L107 -> L207  (cells C8, C9)
#   - ELECTIVE_B4 / ELECTIVE_B5 basket rows (23-32): lecture-slot classrooms
#     are renumbered, and the tutorial-slot cell (column E) now also shows
#     the classroom that the lecture uses.
#
# NOTE: this runtime's PowerShell subset does not bind named (-param)
# arguments on user-defined functions, so helper calls below use
# positional arguments only.

$wb = $excel.ActiveWorkbook

$sectionASheets = @("Regular_Section_A", "PreMid_Section_A", "PostMid_Section_A")
$sectionBSheets = @("Regular_Section_B", "PreMid_Section_B", "PostMid_Section_B")

# Sets the lecture-slot cell (column D) and tutorial-slot cell (column E)
# for one elective-basket row. $d1/$d2 are the two lecture day/time
# strings (same room for both), $e1 is the tutorial day/time string.
function Set-ElectiveRow {
    param($ws, $row, $d1, $d2, $room, $e1)
    $ws.Range("D$row").Value = "$d1 [$room], $d2 [$room]"
    $ws.Range("E$row").Value = "$e1 [$room]"
}

foreach ($name in $sectionASheets) {
    $ws = $wb.Worksheets.Item($name)

    # CS303 (Lab) room: L107 -> L207
    $ws.Range("C8").Value = "CS303 (Lab) [L207]"
    $ws.Range("C9").Value = "CS303 (Lab) [L207]"

    # ELECTIVE_B4 (Tue/Thu lecture, Wed tutorial)
    Set-ElectiveRow $ws 23 "Tue 13:00-14:30" "Thu 13:00-14:30" "C101" "Wed 14:30-15:30"
    Set-ElectiveRow $ws 24 "Tue 13:00-14:30" "Thu 13:00-14:30" "C102" "Wed 14:30-15:30"
    Set-ElectiveRow $ws 25 "Tue 13:00-14:30" "Thu 13:00-14:30" "C104" "Wed 14:30-15:30"
    Set-ElectiveRow $ws 26 "Tue 13:00-14:30" "Thu 13:00-14:30" "C202" "Wed 14:30-15:30"

    # ELECTIVE_B5 (Mon/Wed lecture, Thu tutorial)
    Set-ElectiveRow $ws 27 "Mon 15:30-17:00" "Wed 15:30-17:00" "C101" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 28 "Mon 15:30-17:00" "Wed 15:30-17:00" "C102" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 29 "Mon 15:30-17:00" "Wed 15:30-17:00" "C104" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 30 "Mon 15:30-17:00" "Wed 15:30-17:00" "C202" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 31 "Mon 15:30-17:00" "Wed 15:30-17:00" "C203" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 32 "Mon 15:30-17:00" "Wed 15:30-17:00" "C204" "Thu 14:30-15:30"
}

foreach ($name in $sectionBSheets) {
    $ws = $wb.Worksheets.Item($name)

    # CS303 (Lab) room: L107 -> L207
    $ws.Range("C8").Value = "CS303 (Lab) [L207]"
    $ws.Range("C9").Value = "CS303 (Lab) [L207]"

    # ELECTIVE_B4 (Tue/Thu lecture, Wed tutorial)
    Set-ElectiveRow $ws 23 "Tue 13:00-14:30" "Thu 13:00-14:30" "C101" "Wed 14:30-15:30"
    Set-ElectiveRow $ws 24 "Tue 13:00-14:30" "Thu 13:00-14:30" "C102" "Wed 14:30-15:30"
    Set-ElectiveRow $ws 25 "Tue 13:00-14:30" "Thu 13:00-14:30" "C104" "Wed 14:30-15:30"
    Set-ElectiveRow $ws 26 "Tue 13:00-14:30" "Thu 13:00-14:30" "C202" "Wed 14:30-15:30"

    # ELECTIVE_B5 (Mon/Wed lecture, Thu tutorial)
    Set-ElectiveRow $ws 27 "Mon 15:30-17:00" "Wed 15:30-17:00" "C101" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 28 "Mon 15:30-17:00" "Wed 15:30-17:00" "C102" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 29 "Mon 15:30-17:00" "Wed 15:30-17:00" "C104" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 30 "Mon 15:30-17:00" "Wed 15:30-17:00" "C202" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 31 "Mon 15:30-17:00" "Wed 15:30-17:00" "C203" "Thu 14:30-15:30"
    Set-ElectiveRow $ws 32 "Mon 15:30-17:00" "Wed 15:30-17:00" "C204" "Thu 14:30-15:30"
}

Write-Host "Done applying elective lecture/tutorial scheduling fix."
